# Bugs.xlsx — add bug rows to the tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: B1 should share the exact same (red-fill / bold / white / centered / wrap)
# style as A1 and C1 instead of the near-duplicate style that never actually applied the fill. ---
$ws.Range("B1").Interior.Color = 6711008   # BGR for FFE06666 -> matches A1 / C1 header fill

# --- New bug rows ---
$rows = @(
    @(
        "No damage",
        "Player's hp will never reach zero, unless you were actually able to attack someone with a 0 damage attack (not intended)",
        "Modify the takeDamage method, total/current health should not be equal to damage taken, it should have damage numbers subtracted from it."
    ),
    @(
        "Game cannot end",
        "Due to previous bug, game will never be finished.",
        "See above."
    ),
    @(
        "Cannot avoid damage",
        "There's no option to avoid incoming damage even though output clearly says you can either take or avoid damage",
        "Implement a method that allows for players to roll a dice which allows them to negate or entirely avoid incoming damage."
    )
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# --- Wrap text on the Description/Fix cells of rows 2 & 3 ---
$ws.Range("B2:C2").WrapText = $true
$ws.Range("B3").WrapText = $true

# --- Wrap text + vertically centered on row 4's Description/Fix cells ---
$ws.Range("B4:C4").WrapText = $true
$ws.Range("B4:C4").VerticalAlignment = -4108   # xlCenter

# --- Row heights to fit the wrapped content ---
$ws.Rows.Item(2).RowHeight = 57
$ws.Rows.Item(3).RowHeight = 39.75
$ws.Rows.Item(4).RowHeight = 73.5

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 36.29
$ws.Columns.Item(2).ColumnWidth = 40.29
$ws.Columns.Item(3).ColumnWidth = 41.57

# --- Freeze header row and leave the final selection on D4, matching the authored view ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D4").Select()

Write-Output "done"
